# Fixzone.xlsx update - commit "2019-06-22 kl. 19:28"
# - Fix "Augisti" -> "Augusti" typo in the month label (J20)
# - Mark days 20-25 (rows 9-14) as completed: green fill formatting + 5 hours logged in column I
#   (row 11's D cell keeps its existing highlight and is left untouched)
# - Move the active cell selection to C30

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an already-"done" cell (E4) as the format source and paste its formatting
# (green fill + thin border) onto the newly completed day rows.
$ws.Range("E4").Copy()
$ws.Range("D9:H10").PasteSpecial(-4122)

$ws.Range("E4").Copy()
$ws.Range("E11:H11").PasteSpecial(-4122)

$ws.Range("E4").Copy()
$ws.Range("D12:H14").PasteSpecial(-4122)

# Log 5 hours for each of the newly completed days
$ws.Range("I9").Value = 5
$ws.Range("I10").Value = 5
$ws.Range("I11").Value = 5
$ws.Range("I12").Value = 5
$ws.Range("I13").Value = 5
$ws.Range("I14").Value = 5

# Fix the misspelled month label
$ws.Range("J20").Value = "Augusti"

# Leave the selection where the author last clicked before saving
$ws.Range("C30").Select()
